# Updated testCasesSample.xlsx with new test cases
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCaseTemplate")

# Rename the template title shown in the merged header row (A1)
$ws.Range("A1").Value = "Test Case San Diego Ticketing System"

# Update who executed the GoogleSearch_1 test case
$ws.Range("J3").Value = "Nick"

# Leave the sheet with the same cell selection recorded in the saved file
$ws.Range("M5").Select() | Out-Null
